$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.513.59'
$ws.Range('E2').Value = '  +4.34%  '

$ws.Range('D3').Value = '1.842.16'
$ws.Range('E3').Value = '  +3.77%  '

$ws.Range('D4').Value = '1.030'
$ws.Range('E4').Value = '  +2.91%  '

$ws.Range('D5').Value = '319.25'
$ws.Range('E5').Value = '  +4.33%  '

$ws.Range('D6').Value = '1.026'
$ws.Range('E6').Value = '  +2.54%  '

$ws.Range('D7').Value = '0.4372'
$ws.Range('E7').Value = '  +3.43%  '

$ws.Range('D8').Value = '0.3733'
$ws.Range('E8').Value = '  +3.73%  '

$ws.Range('D9').Value = '0.07386'
$ws.Range('E9').Value = '  +3.60%  '

$ws.Range('D10').Value = '0.8743'
$ws.Range('E10').Value = '  +4.67%  '

$ws.Range('D11').Value = '21.43'
$ws.Range('E11').Value = '  +4.93%  '

$ws.Range('D12').Value = '1.868.21'
$ws.Range('E12').Value = '  +5.08%  '

$ws.Range('D13').Value = '5.478'
$ws.Range('E13').Value = '  +4.58%  '

$ws.Range('D14').Value = '6.700'
$ws.Range('E14').Value = '  +3.84%  '

$ws.Range('D15').Value = '0.07157'
$ws.Range('E15').Value = '  +4.25%  '

$ws.Range('D16').Value = '82.85'
$ws.Range('E16').Value = '  +4.96%  '

$ws.Range('D17').Value = '1.034'
$ws.Range('E17').Value = '  +3.22%  '

$ws.Range('D18').Value = '0.000009016'
$ws.Range('E18').Value = '  +4.43%  '

$ws.Range('D19').Value = '1.027'
$ws.Range('E19').Value = '  +2.61%  '

$ws.Range('D20').Value = '15.43'
$ws.Range('E20').Value = '  +3.54%  '

$ws.Range('D21').Value = '27.536.69'
$ws.Range('E21').Value = '  +4.42%  '

$ws.Range('D22').Value = '5.247'
$ws.Range('E22').Value = '  +3.24%  '

$ws.Range('D23').Value = '11.29'
$ws.Range('E23').Value = '  +3.18%  '

$ws.Range('D24').Value = '2.080.41'
$ws.Range('E24').Value = '  +3.99%  '

$ws.Range('D25').Value = '157.04'
$ws.Range('E25').Value = '  +3.24%  '

$ws.Range('D26').Value = '1.914'
$ws.Range('E26').Value = '  +5.56%  '

$ws.Range('D27').Value = '18.66'
$ws.Range('E27').Value = '  +3.76%  '

$ws.Range('D28').Value = '5.259'
$ws.Range('E28').Value = '  +3.87%  '

$ws.Range('D29').Value = '1.929'
$ws.Range('E29').Value = '  +5.46%  '

$ws.Range('D30').Value = '116.43'
$ws.Range('E30').Value = '  +1.64%  '

$ws.Range('D31').Value = '0.09072'
$ws.Range('E31').Value = '  +2.50%  '

$ws.Range('E32').Value = '  +7.54%  '

$ws.Range('D33').Value = '0.7621'
$ws.Range('E33').Value = '  +4.93%  '

$ws.Range('D34').Value = '4.492'
$ws.Range('E34').Value = '  +4.03%  '

$ws.Range('D35').Value = '2.876'
$ws.Range('E35').Value = '  +5.30%  '

$ws.Range('D36').Value = '1.029'
$ws.Range('E36').Value = '  +2.90%  '

$ws.Range('D37').Value = '1.149'
$ws.Range('E37').Value = '  +5.26%  '

$ws.Range('D38').Value = '0.01969'
$ws.Range('E38').Value = '  +4.55%  '

$ws.Range('D39').Value = '0.05261'
$ws.Range('E39').Value = '  +2.58%  '

$ws.Range('D40').Value = '0.5175'
$ws.Range('E40').Value = '  +5.25%  '

$ws.Range('D41').Value = '2.788'
$ws.Range('E41').Value = '  +7.31%  '

$ws.Range('D42').Value = '0.1666'
$ws.Range('E42').Value = '  +3.62%  '

$ws.Range('D43').Value = '6.554'
$ws.Range('E43').Value = '  +3.60%  '

$ws.Range('D44').Value = '8.505'
$ws.Range('E44').Value = '  +6.81%  '

$ws.Range('D45').Value = '109.24'
$ws.Range('E45').Value = '  +4.41%  '

$ws.Range('D46').Value = '10.57'
$ws.Range('E46').Value = '  +4.23%  '

$ws.Range('D47').Value = '1.031'
$ws.Range('E47').Value = '  +2.98%  '

$ws.Range('E48').Value = '  +3.18%  '

$ws.Range('B49').Value = 'RenderToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D49').Value = '1.911'
$ws.Range('E49').Value = '  +10.58%  '

$ws.Range('B50').Value = 'Decentraland'
$ws.Range('C50').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D50').Value = '0.4642'
$ws.Range('E50').Value = '  +4.55%  '

$ws.Range('D51').Value = '0.06325'
$ws.Range('E51').Value = '  +2.49%  '
